$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72 (this shifts existing rows 72..223 down to 73..224,
# matching the diff which moves every record at old row N to new row N+1).
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new weekly record.
$ws.Range("A72").Value = 5
$ws.Range("B72").Value = "Macroferia Regional de Talca"
$ws.Range("C72").Value = "Maule"
$ws.Range("D72").Value = 44544
$ws.Range("E72").Value = 7
$ws.Range("F72").Value = 100114014
$ws.Range("G72").Value = "Betarraga"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 550
$ws.Range("L72").Value = 550
$ws.Range("M72").Value = 550
$ws.Range("N72").Value = '$/paquete 5 unidades'
$ws.Range("O72").Value = "Región del Maule"
$ws.Range("P72").Value = 110
$ws.Range("Q72").Value = 5
$ws.Range("R72").Value = "Hortaliza"
